$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$vins = @(
    "8LDC2230CE021610",
    "8LDC2230CE021641",
    "KNAPM81ABG7000347",
    "KNALN414BB5015438",
    "KNALN414BC5079465",
    "KNAHU812AH7158037",
    "KNAFT411BD5785277",
    "KNAMC812BG6082657",
    "KNALD225385130724",
    "KNABX512BCT035612",
    "KNALT412BE6010919",
    "KNAMB761286195778",
    "KNAJX81EFF7000871",
    "8LCDC22328E008536"
)

# Clear out the existing column B values (4564, 46, 46, 46456) and the
# numeric values in A2:A4 (3423, 4354, 54654)
$ws.Range("B1:B4").ClearContents()
$ws.Range("A2:A4").ClearContents()

# A3 previously had the default style; the target keeps it bordered like A2,
# so copy A2's format (border, no fill) onto A3 before writing new values.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)

# Write the VIN numbers into A2:A15
for ($i = 0; $i -lt $vins.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $vins[$i]
}

# Select A15 as the active cell, matching the final selection state
$ws.Range("A15").Select()
